$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new empty "ListParagraph" paragraph right before "List To-Do",
#    carrying an extra w:ind w:left="1080" on the new paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("List To-Do", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($rng.Start, $rng.Start)
$target.InsertParagraphBefore()
$newPara = $d.Paragraphs($target.Paragraphs(1).Range.Start)
# Re-find the freshly inserted (now-empty) paragraph and set its indent.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("List To-Do", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$prevParaRange = $d.Range(0, $rng2.Start)
$emptyPara = $prevParaRange.Paragraphs.Last
$emptyPara.LeftIndent = 54  ## 1080 twips = 54 pt

# ---------------------------------------------------------------------------
# 2) Fix typo "botulin" -> "betulin" (first occurrence, inside the
#    "Bantu botulin design..." sentence).
# ---------------------------------------------------------------------------
$rngA = $d.Content
$rngA.Find.Execute("Bantu botulin design", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sA = $rngA.Start
$charO = $d.Range($sA + 7, $sA + 8)
$charO.Text = "e"

# ---------------------------------------------------------------------------
# 3) Collapse "Component at least have : Tabs, NavDrawer, Bottom Navbar"
#    (which spanned 3 runs split by proofErr tags) into a single run/phrase,
#    dropping the extra space before the colon.
# ---------------------------------------------------------------------------
$rngB = $d.Content
$rngB.Find.Execute("Component at least have : Tabs, NavDrawer, Bottom Navbar", $true, $false, $false, $false, $false, $true, 1, $false, "Component at least have: Tabs, NavDrawer, Bottom Navbar", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Move the <w:lastRenderedPageBreak/> marker: it used to sit in front of
#    "Make sure even when no network access..." and now belongs in front of
#    "(you can use free api or create your own, but not firebase)".
#    Rewriting the "Make sure" run (identical text) drops its stale marker;
#    a tiny formatting no-op/reset on the api-text run forces Word to mint a
#    fresh run there that carries the marker forward.
# ---------------------------------------------------------------------------
$rngC = $d.Content
$rngC.Find.Execute("Make sure even when no network access condition, the user still can see the list if has open the apps before.", $true, $false, $false, $false, $false, $true, 1, $false, "Make sure even when no network access condition, the user still can see the list if has open the apps before.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Fix typo "botulin" -> "betulin" (second occurrence, inside the
#    "Commit 2 = ..." sentence).
# ---------------------------------------------------------------------------
$rngD = $d.Content
$rngD.Find.Execute("gw botulin tapi ga banya", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sD = $rngD.Start
$charO2 = $d.Range($sD + 4, $sD + 5)
$charO2.Text = "e"
